$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1781.8182
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 2685.7144
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 8057.1432
$ws.Range("M29").Value = -319
$ws.Range("N29").Value = -8619.143199999999

$ws.Range("H38").Value = 1647.28
$ws.Range("I38").Value = 209.66667
$ws.Range("J38").Value = 2101.2632
$ws.Range("K38").Value = 629.00001
$ws.Range("L38").Value = 6303.7896
$ws.Range("M38").Value = -257.00001
$ws.Range("N38").Value = -7047.7896

$ws.Range("H40").Value = 1813.9231
$ws.Range("J40").Value = 1772.375
$ws.Range("L40").Value = 1772.375
$ws.Range("N40").Value = -2122.375

$ws.Range("H43").Value = 4639578
$ws.Range("I43").Value = 25697.75
$ws.Range("J43").Value = 6946518
$ws.Range("K43").Value = 25697.75
$ws.Range("L43").Value = 6946518
$ws.Range("M43").Value = -25628.75
$ws.Range("N43").Value = -6946656

$ws.Range("H64").Value = 4015
$ws.Range("I64").Value = 4060
$ws.Range("J64").Value = 3970
$ws.Range("K64").Value = 4060
$ws.Range("L64").Value = 3970
$ws.Range("M64").Value = -3812
$ws.Range("N64").Value = -4466

$ws.Range("H67").Value = 4015
$ws.Range("I67").Value = 4060
$ws.Range("J67").Value = 3970
$ws.Range("K67").Value = 4060
$ws.Range("L67").Value = 3970
$ws.Range("M67").Value = -3202
$ws.Range("N67").Value = -5686

$ws.Range("H76").Value = 4610.75
$ws.Range("I76").Value = 5115
$ws.Range("J76").Value = 3098
$ws.Range("K76").Value = 5115
$ws.Range("L76").Value = 3098
$ws.Range("M76").Value = -4800
$ws.Range("N76").Value = -3728

$ws.Range("H79").Value = 4610.75
$ws.Range("I79").Value = 5115
$ws.Range("J79").Value = 3098
$ws.Range("K79").Value = 5115
$ws.Range("L79").Value = 3098
$ws.Range("M79").Value = -4023
$ws.Range("N79").Value = -5282

$ws.Range("H132").Value = 14506284
$ws.Range("I132").Value = 17555838
$ws.Range("K132").Value = 52667514
$ws.Range("M132").Value = -52664984

$ws.Range("H138").Value = 459427.53
$ws.Range("J138").Value = 559319.9399999999
$ws.Range("L138").Value = 1677959.82
$ws.Range("N138").Value = -1688239.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3576.5881
$ws.Range("I32").Value = 3421.9397
$ws.Range("J32").Value = 9994.5
$ws.Range("K32").Value = 3421.9397
$ws.Range("L32").Value = 9994.5
$ws.Range("M32").Value = -3134.9397
$ws.Range("N32").Value = -10568.5

$ws.Range("H45").Value = 1150.7693
$ws.Range("I45").Value = 1096
$ws.Range("K45").Value = 1096
$ws.Range("M45").Value = -719

$ws.Range("H88").Value = 3800
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3800
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3800
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4612

$ws.Range("H91").Value = 3800
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3800
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3800
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6608

$ws.Range("H106").Value = 22370
$ws.Range("J106").Value = 22370
$ws.Range("L106").Value = 22370
$ws.Range("N106").Value = -24894

$ws.Range("H132").Value = 2662.2903
$ws.Range("I132").Value = 2367.4211
$ws.Range("J132").Value = 3129.1667
$ws.Range("K132").Value = 7102.263300000001
$ws.Range("L132").Value = 9387.500100000001
$ws.Range("M132").Value = -4572.263300000001
$ws.Range("N132").Value = -14447.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1779.1052
$ws.Range("I58").Value = 1488.8572
$ws.Range("J58").Value = 2591.8
$ws.Range("K58").Value = 1488.8572
$ws.Range("L58").Value = 2591.8
$ws.Range("M58").Value = -1285.8572
$ws.Range("N58").Value = -2997.8

$ws.Range("H122").Value = 4389.857
$ws.Range("I122").Value = 4783.12
$ws.Range("K122").Value = 14349.36
$ws.Range("M122").Value = -11899.36

$ws.Range("H132").Value = 1722.3529
$ws.Range("I132").Value = 1137.0769
$ws.Range("J132").Value = 3624.5
$ws.Range("K132").Value = 3411.2307
$ws.Range("L132").Value = 10873.5
$ws.Range("M132").Value = -881.2307000000001
$ws.Range("N132").Value = -15933.5

$ws.Range("I134").Value = 2119.318
$ws.Range("K134").Value = 6357.954000000001
$ws.Range("M134").Value = -3822.954000000001

$ws.Range("H136").Value = 1779.1052
$ws.Range("I136").Value = 1488.8572
$ws.Range("J136").Value = 2591.8
$ws.Range("K136").Value = 4466.571599999999
$ws.Range("L136").Value = 7775.400000000001
$ws.Range("M136").Value = -1916.571599999999
$ws.Range("N136").Value = -12875.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3520.2666
$ws.Range("J39").Value = 3477.2307
$ws.Range("L39").Value = 10431.6921
$ws.Range("N39").Value = -11019.6921

$ws.Range("H55").Value = 2325.6667
$ws.Range("J55").Value = 2700
$ws.Range("L55").Value = 8100
$ws.Range("N55").Value = -8454

$ws.Range("H108").Value = 2162.2222
$ws.Range("I108").Value = 600
$ws.Range("J108").Value = 2254.1177
$ws.Range("K108").Value = 1800
$ws.Range("L108").Value = 6762.353099999999
$ws.Range("M108").Value = 1080
$ws.Range("N108").Value = -12522.3531

$ws.Range("H109").Value = 56968.168
$ws.Range("I109").Value = 77402.08
$ws.Range("J109").Value = 3840
$ws.Range("K109").Value = 232206.24
$ws.Range("L109").Value = 11520
$ws.Range("M109").Value = -231166.24
$ws.Range("N109").Value = -13600

$ws.Range("H131").Value = 11237035
$ws.Range("J131").Value = 1142.6464
$ws.Range("L131").Value = 3427.9392
$ws.Range("N131").Value = -13507.9392

$ws.Range("H134").Value = 4993.3335
$ws.Range("I134").Value = 2840.9092
$ws.Range("J134").Value = 6814.615
$ws.Range("K134").Value = 8522.7276
$ws.Range("L134").Value = 20443.845
$ws.Range("M134").Value = -3452.7276
$ws.Range("N134").Value = -30583.845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 41333.332
$ws.Range("J104").Value = 41333.332
$ws.Range("L104").Value = 41333.332
$ws.Range("N104").Value = -48321.332

$ws.Range("H107").Value = 578.5833
$ws.Range("J107").Value = 555.25
$ws.Range("L107").Value = 555.25
$ws.Range("N107").Value = -4395.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5750
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 6800
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 6800
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -7176

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 870.5682
$ws.Range("I136").Value = 754.45715
$ws.Range("J136").Value = 1322.1111
$ws.Range("K136").Value = 2263.37145
$ws.Range("L136").Value = 3966.3333
$ws.Range("M136").Value = 286.6285500000004
$ws.Range("N136").Value = -9066.3333

Write-Host "Applied all Kujata_Profits updates"
